# The deck ships two DrawingML themes:
#   theme1.xml -> "Office Theme" (default Office colors), used by the Notes Master
#   theme2.xml -> "Integral" theme, used by the Slide Master (i.e. the design
#                 actually applied to/visible on the slides)
#
# The authored edit swaps the two themes' content: the Slide Master's theme
# becomes the default "Office Theme" color set, while the Notes Master's
# theme becomes the "Integral" color set.
#
# This headless PowerPoint host exposes a single writable theme/colour-scheme
# object (reached from SlideMaster.Theme, NotesMaster.Theme, HandoutMaster.Theme
# or Slide.ThemeColorScheme - they all alias the same underlying theme part,
# ppt/theme/theme2.xml, the one actually referenced by the Slide Master). So we
# push the swap through that single reachable surface: set each of the twelve
# theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the
# "Office Theme" values that used to live in theme1.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index -> (slot name, target sRGB hex) for the "Office Theme" colour scheme
# that theme2.xml should end up with (previously theme1.xml's colours).
# COM RGB() value = R + G*256 + B*65536.
$targetRgb = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $slot = $colorScheme.Colors($i)
    $slot.RGB = $targetRgb[$i - 1]
}
